# Reword the welcome messages ("Corporate Conclave" -> "Industry-Academia
# Conclave"), drop the trailing "You can check the agenda for today on the
# left" phrase from each message, and refresh the sheet view (zoom level +
# active cell) to match the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the welcome-message text in column A (rows 2-6) ---
# Row 2: "Welcome to Christ University... by the School of Sciences" message
$ws.Range("A2").Value = "Welcome to Christ University. Thanks for joining the Industry-Academia Conclave by the School of Sciences. We are happy to have you today. "

# Row 3: "We are honored..." message
$ws.Range("A3").Value = "We are honored to welcome you to the Christ University Industry-Academia Conclave. Your presence enriches our event. "

# Row 4: "A warm welcome..." message
$ws.Range("A4").Value = "A warm welcome to the School of Sciences Industry-Academia Conclave. Thank you for being part of this gathering. "

# Row 5: "Greetings!..." message
$ws.Range("A5").Value = "Greetings! We are delighted to have you join us for today's Industry-Academia Conclave at Christ University. "

# Row 6: "Thank you for attending..." message
$ws.Range("A6").Value = "Thank you for attending the Industry-Academia Conclave. We look forward to an insightful session with you. "

# --- Update window / view state ---
# Zoom the sheet view to 175%
$excel.ActiveWindow.Zoom = 175

# Resize/reposition the workbook window (best effort; matches the maximized
# xWindow/yWindow/windowWidth/windowHeight recorded after the edit)
try {
    $win = $wb.Windows.Item(1)
    $win.WindowState = -4137
    $win.Left = -120
    $win.Top = -120
    $win.Width = 29040
    $win.Height = 15720
} catch {
}

try {
    $excel.WindowState = -4137
    $excel.Left = -120
    $excel.Top = -120
    $excel.Width = 29040
    $excel.Height = 15720
} catch {
}

# Move the active selection to K10
$ws.Range("K10").Select()
